# ModelComponentClassDiagram.pptx edit
# ------------------------------------
# 1) Bump the auto-updating "datetimeFigureOut" footer date fields from
#    4/16/2018 -> 4/17/2018 (slide master + every slide layout).
# 2) Rename the "UndoRedoStack" rectangle on slide 1 to "UndoRedoCareTaker"
#    and widen/shift it so the longer label still fits.

$p = $ppt.ActivePresentation

# --- 1. Date placeholders -------------------------------------------------
# Setting TextRange.Text on a placeholder that currently holds an
# auto-updating field (<a:fld type="datetimeFigureOut">) replaces that
# field with literal text - this mirrors real PowerPoint automation
# behaviour (fields are dynamic; assigning .Text bakes in a static value).

# Slide master's own "Date Placeholder" shape.
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "4/17/2018"
    }
}

# Every slide layout's "Date Placeholder" shape.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $sh = $layout.Shapes.Item($si)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "4/17/2018"
        }
    }
}

# --- 2. UndoRedoStack -> UndoRedoCareTaker -------------------------------
$slide = $p.Slides.Item(1)
for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
    $sh = $slide.Shapes.Item($k)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "UndoRedoStack") {
            $sh.Left = 96
            $sh.Width = 109.2403937
            $sh.TextFrame.TextRange.Text = "UndoRedoCareTaker"
        }
    }
}
